$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.726.83'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.889.35'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '248.32'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4733'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2929'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06538'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.10'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07799'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '96.86'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.886.75'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7387'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.257'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '284.80'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.714.18'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007548'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.135.73'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.326'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.254'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.231'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.98'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.920'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  -1.99%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09750'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.301'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.192'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04865'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.127'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6970'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.723'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01897'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.805'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.348'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '76.10'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.55%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.004'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4282'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8350'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.68'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.512'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.055'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '35.61'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '916.24'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05752'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.89%  '
